$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for all data rows (2-299)
$ws.Range("C2:C299").Value = 45175

# Row 299 gets an explicit row height (matches the committed change)
$ws.Rows.Item(299).RowHeight = 15

# Bring over the number-format / wrap-text styles used on row 299 so the
# new row matches the existing look (date columns B & C, wrapped column R)
$ws.Range("B299").Copy() | Out-Null
$ws.Range("B300").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C299").Copy() | Out-Null
$ws.Range("C300").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("R299").Copy() | Out-Null
$ws.Range("R300").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Add the new record as row 300
$ws.Range("A300").Value = "A 40954-2023"
$ws.Range("B300").Value = 45173
$ws.Range("C300").Value = 45175
$ws.Range("D300").Value = "DALARNAS LÄN"
$ws.Range("E300").Value = "BORLÄNGE"
$ws.Range("G300").Value = 0.7
$ws.Range("H300").Value = 0
$ws.Range("I300").Value = 0
$ws.Range("J300").Value = 0
$ws.Range("K300").Value = 0
$ws.Range("L300").Value = 0
$ws.Range("M300").Value = 0
$ws.Range("N300").Value = 0
$ws.Range("O300").Value = 0
$ws.Range("P300").Value = 0
$ws.Range("Q300").Value = 0
